# Update countries & provincias Spain
#
# 1) Data refresh for "Estados Unidos" (row 4) and "Japon" (row 34).
# 2) "Mayotte" moves up in the (descending, by total cases) ranking to sit
#    right after "Uruguay" / before "Sudan", with refreshed case numbers.
#    That re-sorts the rows in between (Sudan, Georgia, San Marino, Mali)
#    down by one position, while "Maldivas" and everything below keeps its
#    row number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Estados Unidos (row 4) -------------------------------------------------
$ws.Range("B4").Value = 1187510
$ws.Range("C4").Value = 26736
$ws.Range("E4").Value = 940666
$ws.Range("G4").Value = 1137
$ws.Range("H4").Value = 68581

# --- Japon (row 34) ----------------------------------------------------------
$ws.Range("B34").Value = 14877
$ws.Range("C34").Value = 306
$ws.Range("D34").Value = 3981
$ws.Range("E34").Value = 10409
$ws.Range("F34").Value = 321
$ws.Range("G34").Value = 13
$ws.Range("H34").Value = 487

# --- Insert "Mayotte" ahead of "Sudan" (row 109), pushing the following
#     four countries (Sudan, Georgia, San Marino, Mali) down one row -------
$ws.Rows(109).Insert()

$ws.Range("A109").Value = "Mayotte"
$ws.Range("B109").Value = 650
$ws.Range("C109").Value = 54
$ws.Range("D109").Value = 235
$ws.Range("E109").Value = 409
$ws.Range("F109").Value = 4
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 6

# The old "Mayotte" row (formerly row 113, now pushed down to row 114 by the
# insert above) is now redundant - remove it so the row count, and every
# country from "Maldivas" onward, stays exactly where it was.
$ws.Rows(114).Delete()
